{"js": "// Replace the run(s) containing \"Qwertyui\" with a single run of 40 \"1\"s,\n// leaving the surrounding paragraph (and its bookmark) untouched.\nconst body = context.document.body;\nconst results = body.search(\"Qwertyui\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nresults.items[0].insertText(\n  \"1111111111111111111111111111111111111111\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Replace the \"Qwertyui\" run text with 40 \"1\" characters, leaving the\n# paragraph (and its _GoBack bookmark) otherwise untouched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Qwertyui\"\n$find.Replacement.Text = \"1111111111111111111111111111111111111111\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
